$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 876.2
$ws.Range("I98").Value = 791.7692
$ws.Range("K98").Value = 791.7692
$ws.Range("M98").Value = 706.2308
# Row 122
$ws.Range("H122").Value = 876.2
$ws.Range("I122").Value = 791.7692
$ws.Range("K122").Value = 2375.3076
$ws.Range("M122").Value = 74.69239999999991
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 125
$ws.Range("H125").Value = 31186
$ws.Range("I125").Value = 41081.332
$ws.Range("J125").Value = 1500
$ws.Range("K125").Value = 369731.988
$ws.Range("L125").Value = 13500
$ws.Range("M125").Value = -367271.988
$ws.Range("N125").Value = -18420
# Row 126
$ws.Range("H126").Value = 40000
$ws.Range("J126").Value = 40000
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880
# Row 128
$ws.Range("H128").Value = 31920
$ws.Range("J128").Value = 31920
$ws.Range("L128").Value = 31920
$ws.Range("N128").Value = -41880
# Row 129
$ws.Range("H129").Value = 828.5714
$ws.Range("J129").Value = 980
$ws.Range("L129").Value = 2940
$ws.Range("N129").Value = -12940
# Row 131
$ws.Range("H131").Value = 2245.9722
$ws.Range("I131").Value = 1161.1666
$ws.Range("K131").Value = 3483.4998
$ws.Range("M131").Value = 1556.5002
# Row 132
$ws.Range("H132").Value = 2456.0715
$ws.Range("I132").Value = 2054.3416
$ws.Range("J132").Value = 3554.1333
$ws.Range("K132").Value = 6163.024800000001
$ws.Range("L132").Value = 10662.3999
$ws.Range("M132").Value = -3633.024800000001
$ws.Range("N132").Value = -15722.3999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 3112.7754
$ws.Range("I132").Value = 3374.2563
$ws.Range("J132").Value = 2093
$ws.Range("K132").Value = 10122.7689
$ws.Range("L132").Value = 6279
$ws.Range("M132").Value = -7592.768899999999
$ws.Range("N132").Value = -11339

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 50404.87
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 50404.87
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 50404.87
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -50994.87
# Row 34
$ws.Range("H34").Value = 50404.87
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 50404.87
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 50404.87
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -50808.87
# Row 58
$ws.Range("H58").Value = 1076.409
$ws.Range("I58").Value = 904.7692
$ws.Range("J58").Value = 1324.3334
$ws.Range("K58").Value = 904.7692
$ws.Range("L58").Value = 1324.3334
$ws.Range("M58").Value = -701.7692
$ws.Range("N58").Value = -1730.3334
# Row 107
$ws.Range("H107").Value = 474.27908
$ws.Range("I107").Value = 458.7
$ws.Range("J107").Value = 510.23077
$ws.Range("K107").Value = 458.7
$ws.Range("L107").Value = 510.23077
$ws.Range("M107").Value = 1461.3
$ws.Range("N107").Value = -4350.23077
# Row 134
$ws.Range("H134").Value = 1750.3334
$ws.Range("I134").Value = 1734.8182
$ws.Range("J134").Value = 1781.3636
$ws.Range("K134").Value = 5204.4546
$ws.Range("L134").Value = 5344.0908
$ws.Range("M134").Value = -2669.4546
$ws.Range("N134").Value = -10414.0908
# Row 136
$ws.Range("H136").Value = 1076.409
$ws.Range("I136").Value = 904.7692
$ws.Range("J136").Value = 1324.3334
$ws.Range("K136").Value = 2714.3076
$ws.Range("L136").Value = 3973.0002
$ws.Range("M136").Value = -164.3076000000001
$ws.Range("N136").Value = -9073.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 259.5
$ws.Range("I107").Value = 283.33334
$ws.Range("J107").Value = 253
$ws.Range("K107").Value = 850.0000200000001
$ws.Range("L107").Value = 759
$ws.Range("M107").Value = 1069.99998
$ws.Range("N107").Value = -4599
# Row 122
$ws.Range("H122").Value = 812.5217
$ws.Range("I122").Value = 416
$ws.Range("J122").Value = 922.6667
$ws.Range("K122").Value = 3744
$ws.Range("L122").Value = 8304.0003
$ws.Range("M122").Value = -1294
$ws.Range("N122").Value = -13204.0003
# Row 123
$ws.Range("H123").Value = 2500
$ws.Range("I123").Value = 1000
$ws.Range("J123").Value = 4000
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 12000
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -16900
# Row 124
$ws.Range("H124").Value = 2000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 2000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 6000
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -15820
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 129
$ws.Range("H129").Value = 10578.454
$ws.Range("I129").Value = 1088.3334
$ws.Range("J129").Value = 21966.6
$ws.Range("K129").Value = 3265.0002
$ws.Range("L129").Value = 65899.79999999999
$ws.Range("M129").Value = 1734.9998
$ws.Range("N129").Value = -75899.79999999999
# Row 130
$ws.Range("H130").Value = 4995
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 4995
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 14985
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -25025
# Row 131
$ws.Range("H131").Value = 871.41174
$ws.Range("I131").Value = 605.0909
$ws.Range("J131").Value = 998.7826
$ws.Range("K131").Value = 1815.2727
$ws.Range("L131").Value = 2996.3478
$ws.Range("M131").Value = 3224.7273
$ws.Range("N131").Value = -13076.3478

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1478.4706
$ws.Range("I7").Value = 1252.4
$ws.Range("J7").Value = 1801.4286
$ws.Range("K7").Value = 1252.4
$ws.Range("L7").Value = 1801.4286
$ws.Range("M7").Value = -1140.4
$ws.Range("N7").Value = -2025.4286
# Row 126
$ws.Range("H126").Value = 1478.4706
$ws.Range("I126").Value = 1252.4
$ws.Range("J126").Value = 1801.4286
$ws.Range("K126").Value = 3757.2
$ws.Range("L126").Value = 5404.2858
$ws.Range("M126").Value = -1287.2
$ws.Range("N126").Value = -10344.2858
# Row 128
$ws.Range("H128").Value = 34276.332
$ws.Range("J128").Value = 34276.332
$ws.Range("L128").Value = 34276.332
$ws.Range("N128").Value = -44236.332
# Row 129
$ws.Range("H129").Value = 40952.668
$ws.Range("J129").Value = 40952.668
$ws.Range("L129").Value = 40952.668
$ws.Range("N129").Value = -50952.668
# Row 130
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
# Row 131
$ws.Range("H131").Value = 34326
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 34326
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 34326
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -44406

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 385.92856
$ws.Range("I107").Value = 365.8889
$ws.Range("J107").Value = 422
$ws.Range("K107").Value = 1097.6667
$ws.Range("L107").Value = 1266
$ws.Range("M107").Value = 822.3333
$ws.Range("N107").Value = -5106
# Row 122
$ws.Range("H122").Value = 8334782
$ws.Range("I122").Value = 15385606
$ws.Range("J122").Value = 1990.4546
$ws.Range("K122").Value = 46156818
$ws.Range("L122").Value = 5971.3638
$ws.Range("M122").Value = -46154368
$ws.Range("N122").Value = -10871.3638
# Row 123
$ws.Range("H123").Value = 21114.334
$ws.Range("J123").Value = 21114.334
$ws.Range("L123").Value = 21114.334
$ws.Range("N123").Value = -30914.334
# Row 124
$ws.Range("H124").Value = 30214.5
$ws.Range("J124").Value = 30214.5
$ws.Range("L124").Value = 30214.5
$ws.Range("N124").Value = -40034.5
# Row 126
$ws.Range("H126").Value = 819.9048
$ws.Range("I126").Value = 609.5625
$ws.Range("J126").Value = 1493
$ws.Range("K126").Value = 1828.6875
$ws.Range("L126").Value = 4479
$ws.Range("M126").Value = 641.3125
$ws.Range("N126").Value = -9419
# Row 127
$ws.Range("H127").Value = 48000
$ws.Range("J127").Value = 48000
$ws.Range("L127").Value = 48000
$ws.Range("N127").Value = -57920
# Row 128
$ws.Range("H128").Value = 36857.5
$ws.Range("J128").Value = 36857.5
$ws.Range("L128").Value = 36857.5
$ws.Range("N128").Value = -46817.5
# Row 130
$ws.Range("H130").Value = 36398
$ws.Range("J130").Value = 36398
$ws.Range("L130").Value = 36398
$ws.Range("N130").Value = -46438
# Row 131
$ws.Range("H131").Value = 35857.5
$ws.Range("J131").Value = 35857.5
$ws.Range("L131").Value = 35857.5
$ws.Range("N131").Value = -45937.5
# Row 132
$ws.Range("H132").Value = 2545.0154
$ws.Range("I132").Value = 2603.7368
$ws.Range("J132").Value = 2126.625
$ws.Range("K132").Value = 7811.2104
$ws.Range("L132").Value = 6379.875
$ws.Range("M132").Value = -5281.2104
$ws.Range("N132").Value = -11439.875
